$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 228.42857
$ws.Cells.Item(4, 9).Value = 130.25
$ws.Cells.Item(4, 10).Value = 359.33334
$ws.Cells.Item(4, 11).Value = 130.25
$ws.Cells.Item(4, 12).Value = 359.33334
$ws.Cells.Item(4, 13).Value = -16.25
$ws.Cells.Item(4, 14).Value = -587.33334

# Row 33
$ws.Cells.Item(33, 8).Value = 143.82353
$ws.Cells.Item(33, 9).Value = 140.41667
$ws.Cells.Item(33, 10).Value = 152
$ws.Cells.Item(33, 11).Value = 140.41667
$ws.Cells.Item(33, 12).Value = 152
$ws.Cells.Item(33, 13).Value = 88.58332999999999
$ws.Cells.Item(33, 14).Value = -610

# Row 51
$ws.Cells.Item(51, 8).Value = 1001
$ws.Cells.Item(51, 9).Value = 1001
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 1001
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -517
$ws.Cells.Item(51, 14).ClearContents()

# Row 126
$ws.Cells.Item(126, 8).Value = 50000
$ws.Cells.Item(126, 10).Value = 50000
$ws.Cells.Item(126, 12).Value = 50000
$ws.Cells.Item(126, 14).Value = -59880

# Row 138
$ws.Cells.Item(138, 8).Value = 2829.0547
$ws.Cells.Item(138, 9).Value = 2180.1667
$ws.Cells.Item(138, 10).Value = 3146.8774
$ws.Cells.Item(138, 11).Value = 6540.500100000001
$ws.Cells.Item(138, 12).Value = 9440.6322
$ws.Cells.Item(138, 13).Value = -1400.500100000001
$ws.Cells.Item(138, 14).Value = -19720.6322

# Row 141
$ws.Cells.Item(141, 8).Value = 4523.95
$ws.Cells.Item(141, 9).Value = 2544.5386
$ws.Cells.Item(141, 10).Value = 8200
$ws.Cells.Item(141, 11).Value = 7633.6158
$ws.Cells.Item(141, 12).Value = 24600
$ws.Cells.Item(141, 13).Value = -2453.6158
$ws.Cells.Item(141, 14).Value = -34960

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 7754492
$ws.Cells.Item(61, 9).Value = 11906873
$ws.Cells.Item(61, 10).Value = 3379.8667
$ws.Cells.Item(61, 11).Value = 11906873
$ws.Cells.Item(61, 12).Value = 3379.8667
$ws.Cells.Item(61, 13).Value = -11906661
$ws.Cells.Item(61, 14).Value = -3803.8667

# Row 74
$ws.Cells.Item(74, 8).Value = 13890921
$ws.Cells.Item(74, 9).Value = 1131.6
$ws.Cells.Item(74, 10).Value = 23812200
$ws.Cells.Item(74, 11).Value = 1131.6
$ws.Cells.Item(74, 12).Value = 23812200
$ws.Cells.Item(74, 13).Value = -257.5999999999999
$ws.Cells.Item(74, 14).Value = -23813948

# Row 77
$ws.Cells.Item(77, 8).Value = 13890921
$ws.Cells.Item(77, 9).Value = 1131.6
$ws.Cells.Item(77, 10).Value = 23812200
$ws.Cells.Item(77, 11).Value = 5658
$ws.Cells.Item(77, 12).Value = 119061000
$ws.Cells.Item(77, 13).Value = -1290
$ws.Cells.Item(77, 14).Value = -119069736

# Row 88
$ws.Cells.Item(88, 8).Value = 1703.8462
$ws.Cells.Item(88, 9).Value = 1375
$ws.Cells.Item(88, 10).Value = 2800
$ws.Cells.Item(88, 11).Value = 1375
$ws.Cells.Item(88, 12).Value = 2800
$ws.Cells.Item(88, 13).Value = -969
$ws.Cells.Item(88, 14).Value = -3612

# Row 91
$ws.Cells.Item(91, 8).Value = 1703.8462
$ws.Cells.Item(91, 9).Value = 1375
$ws.Cells.Item(91, 10).Value = 2800
$ws.Cells.Item(91, 11).Value = 1375
$ws.Cells.Item(91, 12).Value = 2800
$ws.Cells.Item(91, 13).Value = 29
$ws.Cells.Item(91, 14).Value = -5608

# Row 132
$ws.Cells.Item(132, 8).Value = 2336871.5
$ws.Cells.Item(132, 9).Value = 5947.048
$ws.Cells.Item(132, 10).Value = 6415989
$ws.Cells.Item(132, 11).Value = 17841.144
$ws.Cells.Item(132, 12).Value = 19247967
$ws.Cells.Item(132, 13).Value = -15311.144
$ws.Cells.Item(132, 14).Value = -19253027

# Row 136
$ws.Cells.Item(136, 8).Value = 7754492
$ws.Cells.Item(136, 9).Value = 11906873
$ws.Cells.Item(136, 10).Value = 3379.8667
$ws.Cells.Item(136, 11).Value = 35720619
$ws.Cells.Item(136, 12).Value = 10139.6001
$ws.Cells.Item(136, 13).Value = -35718069
$ws.Cells.Item(136, 14).Value = -15239.6001

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Cells.Item(64, 8).Value = 554.2857
$ws.Cells.Item(64, 9).Value = 533.2
$ws.Cells.Item(64, 10).Value = 607
$ws.Cells.Item(64, 11).Value = 533.2
$ws.Cells.Item(64, 12).Value = 607
$ws.Cells.Item(64, 13).Value = -308.2
$ws.Cells.Item(64, 14).Value = -1057

# Row 67
$ws.Cells.Item(67, 8).Value = 554.2857
$ws.Cells.Item(67, 9).Value = 533.2
$ws.Cells.Item(67, 10).Value = 607
$ws.Cells.Item(67, 11).Value = 533.2
$ws.Cells.Item(67, 12).Value = 607
$ws.Cells.Item(67, 13).Value = 246.8
$ws.Cells.Item(67, 14).Value = -2167

# Row 86
$ws.Cells.Item(86, 8).Value = 3076.6667
$ws.Cells.Item(86, 9).Value = 2990
$ws.Cells.Item(86, 10).Value = 3163.3333
$ws.Cells.Item(86, 11).Value = 2990
$ws.Cells.Item(86, 12).Value = 3163.3333
$ws.Cells.Item(86, 13).Value = -1867
$ws.Cells.Item(86, 14).Value = -5409.3333

# Row 88
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()

# Row 89
$ws.Cells.Item(89, 8).Value = 3076.6667
$ws.Cells.Item(89, 9).Value = 2990
$ws.Cells.Item(89, 10).Value = 3163.3333
$ws.Cells.Item(89, 11).Value = 14950
$ws.Cells.Item(89, 12).Value = 15816.6665
$ws.Cells.Item(89, 13).Value = -9334
$ws.Cells.Item(89, 14).Value = -27048.6665

# Row 91
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 88
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()

# Row 91
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()

# Row 97
$ws.Cells.Item(97, 8).Value = 29850
$ws.Cells.Item(97, 10).Value = 29800
$ws.Cells.Item(97, 12).Value = 29800
$ws.Cells.Item(97, 14).Value = -31782

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 707.0833
$ws.Cells.Item(5, 9).Value = 398.63635
$ws.Cells.Item(5, 10).Value = 4100
$ws.Cells.Item(5, 11).Value = 1195.90905
$ws.Cells.Item(5, 12).Value = 12300
$ws.Cells.Item(5, 13).Value = -1083.90905
$ws.Cells.Item(5, 14).Value = -12524

# Row 107
$ws.Cells.Item(107, 8).Value = 452.86957
$ws.Cells.Item(107, 9).Value = 287.73334
$ws.Cells.Item(107, 10).Value = 762.5
$ws.Cells.Item(107, 11).Value = 863.20002
$ws.Cells.Item(107, 12).Value = 2287.5
$ws.Cells.Item(107, 13).Value = 1056.79998
$ws.Cells.Item(107, 14).Value = -6127.5

# Row 131
$ws.Cells.Item(131, 8).Value = 4136.049
$ws.Cells.Item(131, 10).Value = 5207.125
$ws.Cells.Item(131, 12).Value = 15621.375
$ws.Cells.Item(131, 14).Value = -25701.375

# Row 132
$ws.Cells.Item(132, 8).Value = 1882.1052
$ws.Cells.Item(132, 9).Value = 1713.75
$ws.Cells.Item(132, 10).Value = 2004.5454
$ws.Cells.Item(132, 11).Value = 15423.75
$ws.Cells.Item(132, 12).Value = 18040.9086
$ws.Cells.Item(132, 13).Value = -12893.75
$ws.Cells.Item(132, 14).Value = -23100.9086

# Row 133
$ws.Cells.Item(133, 8).Value = 15530.105
$ws.Cells.Item(133, 9).Value = 6411.25
$ws.Cells.Item(133, 10).Value = 22162
$ws.Cells.Item(133, 11).Value = 19233.75
$ws.Cells.Item(133, 12).Value = 66486
$ws.Cells.Item(133, 13).Value = -14173.75
$ws.Cells.Item(133, 14).Value = -76606

# Row 134
$ws.Cells.Item(134, 8).Value = 6147.256
$ws.Cells.Item(134, 9).Value = 2509.875
$ws.Cells.Item(134, 10).Value = 8302.741
$ws.Cells.Item(134, 11).Value = 7529.625
$ws.Cells.Item(134, 12).Value = 24908.223
$ws.Cells.Item(134, 13).Value = -2459.625
$ws.Cells.Item(134, 14).Value = -35048.223

# Row 135
$ws.Cells.Item(135, 8).Value = 707.0833
$ws.Cells.Item(135, 9).Value = 398.63635
$ws.Cells.Item(135, 10).Value = 4100
$ws.Cells.Item(135, 11).Value = 3587.72715
$ws.Cells.Item(135, 12).Value = 36900
$ws.Cells.Item(135, 13).Value = -1052.72715
$ws.Cells.Item(135, 14).Value = -41970

# Row 136
$ws.Cells.Item(136, 8).Value = 4133
$ws.Cells.Item(136, 9).Value = 1588.3334
$ws.Cells.Item(136, 10).Value = 7950
$ws.Cells.Item(136, 11).Value = 4765.0002
$ws.Cells.Item(136, 12).Value = 23850
$ws.Cells.Item(136, 13).Value = 334.9997999999996
$ws.Cells.Item(136, 14).Value = -34050

# Row 137
$ws.Cells.Item(137, 8).Value = 7905
$ws.Cells.Item(137, 9).Value = 13994.111
$ws.Cells.Item(137, 10).Value = 3689.4614
$ws.Cells.Item(137, 11).Value = 41982.333
$ws.Cells.Item(137, 12).Value = 11068.3842
$ws.Cells.Item(137, 13).Value = -36882.333
$ws.Cells.Item(137, 14).Value = -21268.3842

# Row 140
$ws.Cells.Item(140, 8).Value = 1887.7778
$ws.Cells.Item(140, 9).Value = 1887.7778
$ws.Cells.Item(140, 11).Value = 5663.3334
$ws.Cells.Item(140, 13).Value = -483.3334000000004

# Row 141
$ws.Cells.Item(141, 8).Value = 5541.552
$ws.Cells.Item(141, 10).Value = 5874.4736
$ws.Cells.Item(141, 12).Value = 17623.4208
$ws.Cells.Item(141, 14).Value = -27983.4208

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 5237.657
$ws.Cells.Item(70, 9).Value = 5186.1377
$ws.Cells.Item(70, 10).Value = 5486.6665
$ws.Cells.Item(70, 11).Value = 5186.1377
$ws.Cells.Item(70, 12).Value = 5486.6665
$ws.Cells.Item(70, 13).Value = -4916.1377
$ws.Cells.Item(70, 14).Value = -6026.6665

# Row 73
$ws.Cells.Item(73, 8).Value = 5237.657
$ws.Cells.Item(73, 9).Value = 5186.1377
$ws.Cells.Item(73, 10).Value = 5486.6665
$ws.Cells.Item(73, 11).Value = 5186.1377
$ws.Cells.Item(73, 12).Value = 5486.6665
$ws.Cells.Item(73, 13).Value = -4250.1377
$ws.Cells.Item(73, 14).Value = -7358.6665

# Row 80
$ws.Cells.Item(80, 8).Value = 673925.8
$ws.Cells.Item(80, 10).Value = 33458.54
$ws.Cells.Item(80, 12).Value = 33458.54
$ws.Cells.Item(80, 14).Value = -35454.54

# Row 83
$ws.Cells.Item(83, 8).Value = 673925.8
$ws.Cells.Item(83, 10).Value = 33458.54
$ws.Cells.Item(83, 12).Value = 167292.7
$ws.Cells.Item(83, 14).Value = -177276.7

# Row 132
$ws.Cells.Item(132, 8).Value = 6667.926
$ws.Cells.Item(132, 9).Value = 8082
$ws.Cells.Item(132, 10).Value = 4264
$ws.Cells.Item(132, 11).Value = 24246
$ws.Cells.Item(132, 12).Value = 12792
$ws.Cells.Item(132, 13).Value = -21716
$ws.Cells.Item(132, 14).Value = -17852

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Cells.Item(82, 8).Value = 1602.0834
$ws.Cells.Item(82, 9).Value = 1302.6666
$ws.Cells.Item(82, 10).Value = 2500.3333
$ws.Cells.Item(82, 11).Value = 1302.6666
$ws.Cells.Item(82, 12).Value = 2500.3333
$ws.Cells.Item(82, 13).Value = -941.6666
$ws.Cells.Item(82, 14).Value = -3222.3333

# Row 85
$ws.Cells.Item(85, 8).Value = 1602.0834
$ws.Cells.Item(85, 9).Value = 1302.6666
$ws.Cells.Item(85, 10).Value = 2500.3333
$ws.Cells.Item(85, 11).Value = 1302.6666
$ws.Cells.Item(85, 12).Value = 2500.3333
$ws.Cells.Item(85, 13).Value = -54.66660000000002
$ws.Cells.Item(85, 14).Value = -4996.3333
